$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Summary": update aggregate stats to reflect the newly closed trade
# (Trade #10, a MarketMaking loss of -0.04).
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1199.96   # Current Capital
$wsSummary.Range("B4").Value = -0.04     # Total P&L $
$wsSummary.Range("B5").Value = -0.08     # Total P&L %
$wsSummary.Range("B6").Value = 10        # Total Trades
$wsSummary.Range("B8").Value = 6         # Losing Trades
$wsSummary.Range("B9").Value = 40        # Win Rate %

# ---------------------------------------------------------------------------
# Sheet "Strategy Status": update the MarketMaking strategy row (row 4)
# ---------------------------------------------------------------------------
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C4").Value = 99.95999999999999   # Capital
$wsStatus.Range("D4").Value = 10                  # Trades
$wsStatus.Range("E4").Value = -0.04               # P&L $
$wsStatus.Range("F4").Value = -0.04               # P&L %
$wsStatus.Range("G4").Value = 40                  # Win Rate %

# ---------------------------------------------------------------------------
# Append the new Trade #10 row (row 11) to both the "All Trades" sheet and
# the "MarketMaking" sheet - they carry the same trade log content.
# ---------------------------------------------------------------------------
function Add-Trade10Row($ws) {
    $ws.Cells.Item(11, 1).Value = 10

    # Date / Time columns must stay text, not get auto-converted to date/time
    # serials by Excel's smart input parsing - use a leading apostrophe to
    # force text for the date-shaped value.
    $ws.Cells.Item(11, 2).Value = "'2026-02-17"
    $ws.Cells.Item(11, 3).Value = "07:53:13"

    $ws.Cells.Item(11, 4).Value = "MarketMaking"
    $ws.Cells.Item(11, 5).Value = "DOWN"
    $ws.Cells.Item(11, 6).Value = 0.23
    $ws.Cells.Item(11, 7).Value = 0.19
    $ws.Cells.Item(11, 8).Value = "CLOSED"
    $ws.Cells.Item(11, 9).Value = -17.3913
    $ws.Cells.Item(11, 10).Value = -0.04
    $ws.Cells.Item(11, 11).Value = 99.95999999999999
    $ws.Cells.Item(11, 12).Value = 0
    $ws.Cells.Item(11, 13).Value = 0
    $ws.Cells.Item(11, 14).Value = 0.6
    $ws.Cells.Item(11, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item(11, 16).Value = "early_exit"
    $ws.Cells.Item(11, 17).Value = 0.14
}

$wsAllTrades = $wb.Worksheets.Item("All Trades")
Add-Trade10Row $wsAllTrades

$wsMarketMaking = $wb.Worksheets.Item("MarketMaking")
Add-Trade10Row $wsMarketMaking
